# Rename the three header/footer logo pictures.
#   footer1 / footer2 : Pearson Edexcel logo  image1.png -> image2.png
#   header2           : BTec logo             image2.jpg -> image1.jpg
#
# InlineShapes don't expose a settable .Name in the Word object model,
# so each picture is round-tripped through ConvertToShape() (which DOES
# expose .Name), renamed, then converted back to an inline shape.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-FirstPicture($range, [string]$newName) {
    $inline = $range.InlineShapes.Item(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# Footer 1 (odd / default footer): Pearson logo, id="1"
Rename-FirstPicture $sec.Footers.Item(1).Range "image2.png"

# Footer 2 (even-page footer): Pearson logo, id="2"
Rename-FirstPicture $sec.Footers.Item(2).Range "image2.png"

# Header 2 (even-page header): BTec logo, id="3"
Rename-FirstPicture $sec.Headers.Item(2).Range "image1.jpg"
